$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill down column A (Index) for every sector block so each of the 18 rows in
# the block carries the same sector label as the block's header row, instead
# of only the first row of the block.
$ws.Range("A2:A19").Value = "NMI"
$ws.Range("A20:A37").Value = "Business Activity"
$ws.Range("A38:A55").Value = "New Order"
$ws.Range("A56:A73").Value = "Employment"
$ws.Range("A74:A91").Value = "Supplier Deliveries"
$ws.Range("A92:A109").Value = "Inventories"
$ws.Range("A110:A127").Value = "Prices"
$ws.Range("A128:A145").Value = "Backlog Orders"
$ws.Range("A146:A163").Value = "New Export Orders"
$ws.Range("A164:A181").Value = "Imports"
$ws.Range("A182:A199").Value = "Inventory Sentiment"

# Header row: rename "Index Names" -> "Index", "Industries" -> "Industry".
# ("Translation" in C1 is untouched content-wise.)
$ws.Range("A1").Value = "Index"
$ws.Range("B1").Value = "Industry"

# Move the active selection to B1 (matches the saved view state).
$ws.Range("B1").Select()
